# Rename the "units_sheet" worksheet to "units_to_query" and make it the
# active/selected tab (the workbook previously opened with "data_sheet"
# selected; the edit moves the selection to the renamed units sheet).

$wb = $excel.ActiveWorkbook

$unitsSheet = $wb.Worksheets.Item("units_sheet")
$unitsSheet.Name = "units_to_query"

# Activating the sheet sets it as the workbook's active tab (tabSelected on
# the sheet's sheetView / activeTab on the workbook's bookViews), moving the
# "selected" state off of data_sheet and onto units_to_query.
$unitsSheet.Activate()
